# Update the attendance sheet: change date separators from "/" to "-"
# and update the Total/Real/Invalid/Absent counts for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ Date; D; E; G; H }
$rows = @{
    3  = @{ Date = "28-07-2022"; D = 1; E = 0; G = 1; H = 1 }
    4  = @{ Date = "01-08-2022"; D = 1; E = 1; G = 0; H = 0 }
    5  = @{ Date = "04-08-2022"; D = 1; E = 1; G = 0; H = 0 }
    6  = @{ Date = "08-08-2022"; D = 1; E = 1; G = 0; H = 0 }
    7  = @{ Date = "11-08-2022"; D = 1; E = 1; G = 0; H = 0 }
    8  = @{ Date = "15-08-2022"; D = 0; E = 0; G = 0; H = 1 }
    9  = @{ Date = "18-08-2022"; D = 0; E = 0; G = 0; H = 1 }
    10 = @{ Date = "22-08-2022"; D = 1; E = 1; G = 0; H = 0 }
    11 = @{ Date = "25-08-2022"; D = 1; E = 1; G = 0; H = 0 }
    12 = @{ Date = "29-08-2022"; D = 1; E = 1; G = 0; H = 0 }
    13 = @{ Date = "01-09-2022"; D = 1; E = 1; G = 0; H = 0 }
    14 = @{ Date = "05-09-2022"; D = 1; E = 1; G = 0; H = 0 }
    15 = @{ Date = "08-09-2022"; D = 1; E = 1; G = 0; H = 0 }
    16 = @{ Date = "12-09-2022"; D = 0; E = 0; G = 0; H = 1 }
    17 = @{ Date = "15-09-2022"; D = 0; E = 0; G = 0; H = 1 }
    18 = @{ Date = "19-09-2022"; D = 0; E = 0; G = 0; H = 1 }
    19 = @{ Date = "22-09-2022"; D = 0; E = 0; G = 0; H = 1 }
    20 = @{ Date = "26-09-2022"; D = 0; E = 0; G = 0; H = 1 }
    21 = @{ Date = "29-09-2022"; D = 0; E = 0; G = 0; H = 1 }
}

foreach ($r in $rows.Keys) {
    $info = $rows[$r]
    $cellA = $ws.Cells.Item($r, 1)
    # Force the cell to stay text so the dd-mm-yyyy string isn't
    # auto-converted into a date serial number by Excel.
    $cellA.NumberFormat = "@"
    $cellA.Value = $info.Date
    $ws.Cells.Item($r, 4).Value = $info.D
    $ws.Cells.Item($r, 5).Value = $info.E
    $ws.Cells.Item($r, 7).Value = $info.G
    $ws.Cells.Item($r, 8).Value = $info.H
}
